{"js": "// Resume skills update:\n//   \" JavaScript \" -> \" Typescript \" (within the Python/Java/C++/... skills line)\n//   \" AutoCAD \"    -> \" FPGA \"       (within the Arduino/AutoCAD/Lua/... skills line)\n//   \"Lua \"         -> \"Linux \"       (within the Arduino/AutoCAD/Lua/... skills line)\n//   \" Soldering\"   -> \" Verilog\"     (within the Git/SolidWorks/.../Soldering/... skills line)\n//   \" SpatialAnalyzer\" -> \" Systems Engineering\" (same skills line, after the \"|\")\n//\n// Each target word is unique in the body except \"JavaScript\", which also shows up\n// in \"JavaScript (Basic) Certificate\" and a later project blurb -- we only want the\n// very first occurrence (the one inside the skills table). For words that share a\n// paragraph, replacements are issued right-to-left (later word in the line first)\n// so that one edit never disturbs the text another edit still needs to locate.\n\nasync function replaceOnce(context, searchText, replacement, matchCase) {\n  const results = context.document.body.search(searchText, {\n    matchCase: matchCase !== false,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find text: \" + searchText);\n  }\n  results.items[0].insertText(replacement, \"Replace\");\n  await context.sync();\n}\n\n// 1) Skills line: Python | Java | C++ | JavaScript | React | Next.JS | MATLAB |\n//    Only the skill-list \"JavaScript\" (first hit in the document) is touched.\nawait replaceOnce(context, \"JavaScript\", \"Typescript\", true);\n\n// 2) Skills line: Arduino | AutoCAD | Lua | Tailwind CSS | SQL | AngularJS |\n//    \"Lua\" appears after \"AutoCAD\" in the line, so replace it first.\nawait replaceOnce(context, \"Lua\", \"Linux\", true);\nawait replaceOnce(context, \"AutoCAD\", \"FPGA\", true);\n\n// 3) Skills line: Git | SolidWorks | 3D Printing | Soldering | SpatialAnalyzer\n//    \"SpatialAnalyzer\" appears after \"Soldering\" in the line, so replace it first.\nawait replaceOnce(context, \"SpatialAnalyzer\", \"Systems Engineering\", true);\nawait replaceOnce(context, \"Soldering\", \"Verilog\", true);\n", "ps1": "# Resume skills update:\n#   \" JavaScript \" -> \" Typescript \" (Python/Java/C++/... skills line)\n#   \" AutoCAD \"    -> \" FPGA \"       (Arduino/AutoCAD/Lua/... skills line)\n#   \"Lua \"         -> \"Linux \"       (Arduino/AutoCAD/Lua/... skills line)\n#   \" Soldering\"   -> \" Verilog\"     (Git/SolidWorks/.../Soldering/... skills line)\n#   \" SpatialAnalyzer\" -> \" Systems Engineering\" (same skills line, after the \"|\")\n#\n# \"JavaScript\" also appears in \"JavaScript (Basic) Certificate\" and later in a\n# project blurb, so the replace is scoped to the specific skills-line paragraph\n# and limited to a single (first) match (wdReplaceOne) so those other spots are\n# left alone. AutoCAD/Lua/Soldering/SpatialAnalyzer are unique in the document,\n# but are still scoped to their paragraph for safety. Where a paragraph has two\n# replacements, the one that appears later in the line is done first so it can\n# still be located correctly before the earlier edit changes the text before it.\n\n$d = $word.ActiveDocument\n\nfunction Find-ParagraphContaining($doc, $needle) {\n    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {\n        $p = $doc.Paragraphs.Item($i)\n        if ($p.Range.Text.Contains($needle)) {\n            return $p\n        }\n    }\n    return $null\n}\n\nfunction Replace-FirstInRange($rng, $searchText, $replaceText) {\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $searchText\n    $find.Replacement.Text = $replaceText\n    # MatchCase:=True, MatchWholeWord:=True, Wrap:=wdFindStop(0), Replace:=wdReplaceOne(1)\n    $find.Execute($searchText, $true, $true, $false, $false, $false, $true, 0, $false, $replaceText, 1) | Out-Null\n}\n\n# Line: Python | Java | C++ | JavaScript | React | Next.JS | MATLAB |\n$p1 = Find-ParagraphContaining $d \"Python | Java\"\nReplace-FirstInRange $p1.Range \"JavaScript\" \"Typescript\"\n\n# Line: Arduino | AutoCAD | Lua | Tailwind CSS | SQL | AngularJS |\n$p2 = Find-ParagraphContaining $d \"Arduino\"\nReplace-FirstInRange $p2.Range \"Lua\" \"Linux\"\nReplace-FirstInRange $p2.Range \"AutoCAD\" \"FPGA\"\n\n# Line: Git | SolidWorks | 3D Printing | Soldering | SpatialAnalyzer\n$p3 = Find-ParagraphContaining $d \"SolidWorks\"\nReplace-FirstInRange $p3.Range \"SpatialAnalyzer\" \"Systems Engineering\"\nReplace-FirstInRange $p3.Range \"Soldering\" \"Verilog\"\n"}
